# Delete the two rows corresponding to account 004239624 (NINA) and
# account 005135281 (RAFAEL), as removed in the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

$accountsToDelete = @("004239624", "005135281")

for ($r = $rowCount; $r -ge 1; $r--) {
    $valStr = $ws.Cells.Item($r, 1).Text
    if ($accountsToDelete -contains $valStr) {
        $ws.Rows.Item($r).Delete()
    }
}
